$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename "cadastro" -> "cadastroSucesso"
# ---------------------------------------------------------------------------
$sucesso = $wb.Worksheets.Item("cadastro")
$sucesso.Name = "cadastroSucesso"

# ---------------------------------------------------------------------------
# 2) Create "CadastroFalha" as a copy of cadastroSucesso, placed right after
#    it (and before "buscar"), so the resulting sheet order/ids match.
# ---------------------------------------------------------------------------
$sucesso.Copy($null, $sucesso) | Out-Null
$falha = $wb.Worksheets.Item(2)
$falha.Name = "CadastroFalha"

# ---------------------------------------------------------------------------
# 3) Update data on "CadastroFalha"
#    Row2: A=BRUN227 (kept), C=felipe, G=almeida, M=cleared
#    Row3: A=BRUN229,        C=felipe, G=almeida, M=cleared
# ---------------------------------------------------------------------------
$falha.Range("A2").Value = "BRUN227"
$falha.Range("A3").Value = "BRUN229"
$falha.Range("C2").Value = "felipe"
$falha.Range("G2").Value = "almeida"
$falha.Range("C3").Value = "felipe"
$falha.Range("G3").Value = "almeida"
$falha.Range("M2").ClearContents() | Out-Null
$falha.Range("M3").ClearContents() | Out-Null

# Selection on CadastroFalha ends on G3
$falha.Activate() | Out-Null
$falha.Range("G3").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4) Update data on "cadastroSucesso"
#    Row2: A=BRUN230, M=cleared
#    Row3: A=BRUN231, M=cleared
# ---------------------------------------------------------------------------
$sucesso.Range("A2").Value = "BRUN230"
$sucesso.Range("M2").ClearContents() | Out-Null

$sucesso.Range("A3").Value = "BRUN231"
$sucesso.Range("M3").ClearContents() | Out-Null

# Selection on cadastroSucesso stays at A3
$sucesso.Activate() | Out-Null
$sucesso.Range("A3").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5) Update data on "buscar"
# ---------------------------------------------------------------------------
$buscar = $wb.Worksheets.Item("buscar")
$buscar.Range("A1").Value = "Modelo"
$buscar.Range("A2").Value = "HP CHROMEBOOK 14 G1 (ENERGY STAR)"
$buscar.Range("A3").Value = "HP CHROMEBOOK 14 G1 (ES)"
$buscar.Range("A4").Value = "HP ENVY - PORTÁTIL TOQUE 17T"

# ---------------------------------------------------------------------------
# 6) Make "cadastroSucesso" the active tab (tabSelected) as the last action
# ---------------------------------------------------------------------------
$sucesso.Activate() | Out-Null
